$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-17.78504371477973, -0.3018183719541265, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 0.5409336303676118, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, 0.007041171972861475, 0.13761617473549, -17.78504371477973, 4.064919103940984, -17.78504371477973, 2.605558759152097, -17.78504371477973, 3.089088903291963, -17.78504371477973),
  @(-17.78504371477973, -0.1050070424914866, -17.78504371477973, -17.78504371477973, -17.78504371477973, 3.632333430141029, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(3.723424728366882, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(2.763654981791791, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 0.3475614418006542, -17.78504371477973, 1.806414714800739),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, 4.32192202658202, -17.78504371477973, 1.346152602730543, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.191858780280424),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.075985622850666, 3.024458310864056),
  @(-17.78504371477973, -17.78504371477973, 0.4155148990932051, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.648307403053511),
  @(-17.78504371477973, -17.78504371477973, 0.366711151950173, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.796004199543381, -17.78504371477973),
  @(-17.78504371477973, 0.2724456924810373, -0.02275901143399031, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.770842897630567, 0.2278699364096726, 1.498426772281023, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.813374139447759, -0.2996153151975315, 1.624929303632549, -17.78504371477973),
  @(-17.78504371477973, -17.78504371477973, 1.893716300774372, -17.78504371477973, -17.78504371477973, -17.78504371477973, 1.765967916084939, 1.196158410082861, -17.78504371477973, -17.78504371477973),
  @(-17.78504371477973, 3.717600193769315, 3.533115655991288, -17.78504371477973, 1.706419929692664, -17.78504371477973, 0.8551482549862849, 3.700190577250516, -17.78504371477973, 1.558907517784105),
  @(-17.78504371477973, 1.531372621326216, -17.78504371477973, -17.78504371477973, -17.78504371477973, 2.338290883117403, 0.8290657920141765, -17.78504371477973, -17.78504371477973, -17.78504371477973)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}
